# South Dakota_Converted.xlsx update ("Updated policies and graphs"):
#  - Q7: 1 -> 0, T7: 13 -> 12 (denominator used for the T-column percentages)
#  - T24:T221 rescaled from /13 to /12 (values are literal, not formulas)
#  - 12 new daily rows appended (222-233) for 9/30/2020 .. 10/11/2020,
#    following the same pattern as the existing "Dates" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Converted Data")

# --- Row 7 updates -------------------------------------------------------
$ws.Range("Q7").Value = 0
$ws.Range("T7").Value = 12

# --- Rescale the T column for the existing date rows (24-221) ------------
# Values were baked in as count/13; the denominator moved to 12, so every
# existing T value is rescaled in place by 13/12 (this exactly reproduces
# the floating point values in the canonical file, e.g. 0.1538461538461539
# -> 0.1666666666666667).
for ($r = 24; $r -le 221; $r++) {
    $cell = $ws.Cells.Item($r, 20)
    $cell.Value = $cell.Value2 * (13 / 12)
}

# --- Append 12 new daily rows (222-233) for 9/30/2020 .. 10/11/2020 ------
$dates = @(
    "9/30/2020", "10/1/2020", "10/2/2020", "10/3/2020", "10/4/2020",
    "10/5/2020", "10/6/2020", "10/7/2020", "10/8/2020", "10/9/2020",
    "10/10/2020", "10/11/2020"
)

$row = 222
foreach ($d in $dates) {
    # Match the formatting (border/font/alignment) used by the rest of the
    # "Dates" column, then write the date as literal text (leading "'" so
    # it isn't auto-converted to a date serial).
    $ws.Range("A221").Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Cells.Item($row, 1).Value = "'" + $d

    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 1
    $ws.Cells.Item($row, 5).Value = 0
    $ws.Cells.Item($row, 6).Value = 0
    $ws.Cells.Item($row, 7).Value = 0
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0
    $ws.Cells.Item($row, 15).Value = 0
    $ws.Cells.Item($row, 16).Value = 0
    $ws.Cells.Item($row, 17).Value = 0
    $ws.Cells.Item($row, 18).Value = 0
    $ws.Cells.Item($row, 19).Value = 0
    $ws.Cells.Item($row, 20).Value = 1 / 12

    $row++
}
